# Actualizar 02-06-2021 05-09-37
# Appends one new "ping cycle" (14 rows) to the availability log sheet and
# refreshes the timestamp of the previous cycle (rows 1010-1023).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the timestamp of the previous cycle (rows 1010-1023): the
#    "Fecha" (D column) value changes slightly while everything else
#    (Nombre/URL/Disponibilidad) stays identical.
# ---------------------------------------------------------------------
$updatedTimestamp = 44233.19375778936
for ($r = 1010; $r -le 1023; $r++) {
    $ws.Cells.Item($r, 4).Value = $updatedTimestamp
}

# ---------------------------------------------------------------------
# 2) Append the new cycle: rows 1024-1037.
#    Columns: A=Nombre, B=URL (hyperlinked), C=Disponibilidad, D=Fecha
# ---------------------------------------------------------------------
$newTimestamp = 44233.21495916338

$entries = @(
    @{ Name = "Odoo";              Display = "https://www.dataintelligence-group.com/";                 Target = "https://www.dataintelligence-group.com/";                 Loc = "" },
    @{ Name = "Blackbox";          Display = "https://serviciodashboard.azurewebsites.net/";              Target = "https://serviciodashboard.azurewebsites.net/";              Loc = "" },
    @{ Name = "PowerBI";           Display = "https://powerbi.microsoft.com/es-es/";                      Target = "https://powerbi.microsoft.com/es-es/";                      Loc = "" },
    @{ Name = "Dropbox";           Display = "https://www.dropbox.com/";                                  Target = "https://www.dropbox.com/";                                  Loc = "" },
    @{ Name = "Odoo";              Display = "https://dataintelligence.store/";                           Target = "https://dataintelligence.store/";                           Loc = "" },
    @{ Name = "GEE";               Display = "https://app-data-i.users.earthengine.app/";                 Target = "https://app-data-i.users.earthengine.app/";                 Loc = "" },
    @{ Name = "UtilidadesOdoo";    Display = "https://odooutil.azurewebsites.net/";                       Target = "https://odooutil.azurewebsites.net/";                       Loc = "" },
    @{ Name = "Filtros Dashboard"; Display = "https://filtradordashboard.azurewebsites.net/";              Target = "https://filtradordashboard.azurewebsites.net/";              Loc = "" },
    @{ Name = "MapStore";          Display = "https://ide.dataintelligence-group.com/mapstore/#/";         Target = "https://ide.dataintelligence-group.com/mapstore/";          Loc = "/" },
    @{ Name = "GeoServer";         Display = "https://ide.dataintelligence-group.com/geoserver/web/?0";    Target = "https://ide.dataintelligence-group.com/geoserver/web/?0";    Loc = "" },
    @{ Name = "Tomcat";            Display = "https://ide.dataintelligence-group.com/";                    Target = "https://ide.dataintelligence-group.com/";                    Loc = "" },
    @{ Name = "Shiny";             Display = "https://rpubs.com/dataintelligence/";                        Target = "https://rpubs.com/dataintelligence/";                        Loc = "" },
    @{ Name = "Github";            Display = "https://github.com/Sud-Austral/";                            Target = "https://github.com/Sud-Austral/";                            Loc = "" },
    @{ Name = "EZ Exporter";       Display = "https://ezexporter.highviewapps.com/exports/export-profile/";Target = "https://ezexporter.highviewapps.com/exports/export-profile/";Loc = "" }
)

$startRow = 1024
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $entry = $entries[$i]

    $ws.Cells.Item($row, 1).Value = $entry.Name
    $ws.Cells.Item($row, 2).Value = $entry.Display
    $ws.Cells.Item($row, 3).Value = "Disponible"
    $ws.Cells.Item($row, 4).Value = $newTimestamp

    # Apply the same visual styles used by the rest of the log.
    $ws.Cells.Item($row, 2).Style = "Hyperlink"
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $linkCell = $ws.Cells.Item($row, 2)
    if ($entry.Loc -ne "") {
        $ws.Hyperlinks.Add($linkCell, $entry.Target, $entry.Loc)
    } else {
        $ws.Hyperlinks.Add($linkCell, $entry.Target)
    }

    # Re-apply the hyperlink style/value in case Hyperlinks.Add touched them.
    $ws.Cells.Item($row, 2).Value = $entry.Display
    $ws.Cells.Item($row, 2).Style = "Hyperlink"
}

Write-Output "Added rows $startRow..$($startRow + $entries.Count - 1) and refreshed previous cycle timestamps."
